$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates -------------------------------------------------
$ws.Range("A2").Value = "Perfromance Comparison between H5Dwrite_multi and H5Dwrite on Wallaby"
$ws.Range("A4").Value = "# of CHUNKED dsets"
$ws.Range("E4").Value = "# of CONTIG dsets"

# --- New blank spacer row above the title ---------------------------------
$ws.Rows.Item(1).RowHeight = 16.5

# --- Selection (cosmetic, mirrors author's last cursor position) ----------
$ws.Range("A44").Select()

# --- Borders: title row (A2:H2) gets a thick box rule, bold title ---------
$ws.Range("A2:H2").Borders.Item(8).LineStyle = 1
$ws.Range("A2:H2").Borders.Item(8).Weight = -4138
$ws.Range("A2:H2").Borders.Item(9).LineStyle = 1
$ws.Range("A2:H2").Borders.Item(9).Weight = -4138
$ws.Range("A2").Font.Bold = $true

# --- Borders: header row for left table (A4:C4) ----------------------------
$ws.Range("A4:C4").Borders.Item(8).LineStyle = 1
$ws.Range("A4:C4").Borders.Item(9).LineStyle = 1
$ws.Range("A4:C9").Borders.Item(7).LineStyle = 1
$ws.Range("A4:C9").Borders.Item(10).LineStyle = 1
$ws.Range("A9:C9").Borders.Item(9).LineStyle = 1

# --- Borders: header row for right table (E4:G4) ----------------------------
$ws.Range("E4:G4").Borders.Item(8).LineStyle = 1
$ws.Range("E4:G4").Borders.Item(9).LineStyle = 1
$ws.Range("E4:G9").Borders.Item(7).LineStyle = 1
$ws.Range("E4:G9").Borders.Item(10).LineStyle = 1
$ws.Range("E9:G9").Borders.Item(9).LineStyle = 1

# --- Chart series now reference the external workbook (row range shifted) --
$chart1 = $ws.ChartObjects(1).Chart
$chart1.SeriesCollection(1).Formula = "=SERIES(""H5Dwrite"",[1]Sheet1!`$A`$111:`$A`$115,[1]Sheet1!`$B`$111:`$B`$115,1)"
$chart1.SeriesCollection(2).Formula = "=SERIES(""H5Dwrite_multi"",[1]Sheet1!`$A`$111:`$A`$115,[1]Sheet1!`$C`$111:`$C`$115,2)"

$chart2 = $ws.ChartObjects(2).Chart
$chart2.SeriesCollection(1).Formula = "=SERIES(""H5Dwrite"",[1]Sheet1!`$E`$111:`$E`$116,[1]Sheet1!`$F`$111:`$F`$115,1)"
$chart2.SeriesCollection(2).Formula = "=SERIES(""H5Dwrite_multi"",[1]Sheet1!`$E`$111:`$E`$116,[1]Sheet1!`$G`$111:`$G`$115,2)"
